# Auto-generated script applying 2022-11-22 violent crime data updates
# Updates column I (2022 totals) across 46 worksheets, 158 cells total
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 6548
$ws.Range("I3").Value = 6832
$ws.Range("I4").Value = 1572
$ws.Range("I5").Value = 637
$ws.Range("I6").Value = 7853
$ws.Range("I7").Value = 23442

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I7").Value = 743
$ws.Range("I8").Value = 1399
$ws.Range("I9").Value = 119
$ws.Range("I10").Value = 168
$ws.Range("I11").Value = 357
$ws.Range("I15").Value = 272
$ws.Range("I18").Value = 180
$ws.Range("I19").Value = 664
$ws.Range("I20").Value = 581
$ws.Range("I22").Value = 66
$ws.Range("I23").Value = 229
$ws.Range("I27").Value = 207
$ws.Range("I29").Value = 1414
$ws.Range("I30").Value = 78
$ws.Range("I33").Value = 1049
$ws.Range("I36").Value = 321
$ws.Range("I37").Value = 738
$ws.Range("I42").Value = 851
$ws.Range("I43").Value = 206
$ws.Range("I46").Value = 50
$ws.Range("I47").Value = 169
$ws.Range("I48").Value = 302
$ws.Range("I51").Value = 280
$ws.Range("I52").Value = 514
$ws.Range("I53").Value = 257
$ws.Range("I54").Value = 474
$ws.Range("I55").Value = 267
$ws.Range("I63").Value = 78
$ws.Range("I65").Value = 538
$ws.Range("I67").Value = 898
$ws.Range("I71").Value = 67
$ws.Range("I77").Value = 140
$ws.Range("I78").Value = 317
$ws.Range("I79").Value = 666
$ws.Range("I83").Value = 509
$ws.Range("I84").Value = 208
$ws.Range("I85").Value = 1054
$ws.Range("I86").Value = 148
$ws.Range("I88").Value = 218
$ws.Range("I90").Value = 303
$ws.Range("I91").Value = 249
$ws.Range("I95").Value = 355
$ws.Range("I96").Value = 261
$ws.Range("I98").Value = 165
$ws.Range("I100").Value = 39
$ws.Range("I101").Value = 23442

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I3").Value = 403
$ws.Range("I6").Value = 271
$ws.Range("I7").Value = 1054

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I6").Value = 147
$ws.Range("I7").Value = 514

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I3").Value = 76
$ws.Range("I7").Value = 357

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I3").Value = 404
$ws.Range("I6").Value = 449
$ws.Range("I7").Value = 1399

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I6").Value = 120
$ws.Range("I7").Value = 257

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I2").Value = 243
$ws.Range("I3").Value = 229
$ws.Range("I6").Value = 199
$ws.Range("I7").Value = 743

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I2").Value = 78
$ws.Range("I7").Value = 261

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("I2").Value = 21
$ws.Range("I7").Value = 78

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I6").Value = 223
$ws.Range("I7").Value = 738

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I5").Value = 27
$ws.Range("I6").Value = 272
$ws.Range("I7").Value = 898

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("I6").Value = 51
$ws.Range("I7").Value = 208

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I2").Value = 180
$ws.Range("I6").Value = 158
$ws.Range("I7").Value = 538

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I6").Value = 113
$ws.Range("I7").Value = 509

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I3").Value = 128
$ws.Range("I7").Value = 355

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I3").Value = 384
$ws.Range("I6").Value = 337
$ws.Range("I7").Value = 1049

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I2").Value = 102
$ws.Range("I6").Value = 229
$ws.Range("I7").Value = 474

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 417
$ws.Range("I6").Value = 392
$ws.Range("I7").Value = 1414

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I6").Value = 213
$ws.Range("I7").Value = 664

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I4").Value = 42
$ws.Range("I7").Value = 302

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I3").Value = 258
$ws.Range("I6").Value = 304
$ws.Range("I7").Value = 851

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("I5").Value = 2
$ws.Range("I7").Value = 168

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I3").Value = 79
$ws.Range("I4").Value = 45
$ws.Range("I7").Value = 317

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I2").Value = 80
$ws.Range("I7").Value = 267

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("I4").Value = 3
$ws.Range("I7").Value = 50

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I3").Value = 80
$ws.Range("I6").Value = 67
$ws.Range("I7").Value = 229

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I2").Value = 78
$ws.Range("I3").Value = 90
$ws.Range("I7").Value = 249

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I3").Value = 217
$ws.Range("I5").Value = 25
$ws.Range("I7").Value = 666

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I2").Value = 161
$ws.Range("I6").Value = 203
$ws.Range("I7").Value = 581

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("I2").Value = 49
$ws.Range("I6").Value = 84
$ws.Range("I7").Value = 180

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I3").Value = 107
$ws.Range("I7").Value = 321

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("I3").Value = 7
$ws.Range("I6").Value = 39

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("I3").Value = 49
$ws.Range("I7").Value = 169

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("I3").Value = 64
$ws.Range("I7").Value = 272

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("I2").Value = 31
$ws.Range("I7").Value = 165

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("I6").Value = 36
$ws.Range("I7").Value = 119

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("I6").Value = 68
$ws.Range("I7").Value = 218

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I2").Value = 56
$ws.Range("I4").Value = 28
$ws.Range("I7").Value = 207

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I4").Value = 70
$ws.Range("I7").Value = 148

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I6").Value = 106
$ws.Range("I7").Value = 303

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I2").Value = 59
$ws.Range("I3").Value = 79
$ws.Range("I7").Value = 280

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I6").Value = 117
$ws.Range("I7").Value = 206

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("I2").Value = 29
$ws.Range("I7").Value = 66

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("I4").Value = 5
$ws.Range("I7").Value = 67

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("I2").Value = 45
$ws.Range("I7").Value = 140
